$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 481
$ws.Range("B2").Value = "18PG"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 9.3028
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 70
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 0.01547973467254321
$ws.Range("K2").Value = 2.13061800116692

$ws.Range("A3").Value = 5419
$ws.Range("B3").Value = "18PG"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 9.3028
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 25
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 0.04372908258619673
$ws.Range("K3").Value = 2.368985998862328

$ws.Range("A4").Value = 19525
$ws.Range("B4").Value = "DOPE"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 9.784000000000001
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 200
$ws.Range("H4").Value = 75
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 0.8884809258550991
$ws.Range("K4").Value = 9.497002052365326

$ws.Range("A5").Value = 14566
$ws.Range("B5").Value = "DOPE"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 9.784000000000001
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 80
$ws.Range("H5").Value = 55
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 0.9764781998063413
$ws.Range("K5").Value = 10.23952304996591

$ws.Range("A6").Value = 19359
$ws.Range("B6").Value = "DOPE"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 9.784000000000001
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 7
$ws.Range("G6").Value = 350
$ws.Range("H6").Value = 55
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 0.8716161479015142
$ws.Range("K6").Value = 9.354697055992977

$ws.Range("A7").Value = 17213
$ws.Range("B7").Value = "DOPE"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 9.784000000000001
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = 55
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 0.9764781998063413
$ws.Range("K7").Value = 10.23952304996591

$ws.Range("A8").Value = 17115
$ws.Range("B8").Value = "DOPE"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 9.784000000000001
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 90
$ws.Range("H8").Value = 70
$ws.Range("I8").Value = 4
$ws.Range("J8").Value = 0.9764781998063413
$ws.Range("K8").Value = 10.23952304996591

$ws.Range("A9").Value = 24965
$ws.Range("B9").Value = "DOPE"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 9.784000000000001
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 50
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 20
$ws.Range("J9").Value = 0.02657352221973434
$ws.Range("K9").Value = 2.224227380490118

$ws.Range("A10").Value = 27736
$ws.Range("B10").Value = "DOPE"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 9.784000000000001
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 40
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 60
$ws.Range("J10").Value = 0.003856706240723889
$ws.Range("K10").Value = 2.032542887259228

$ws.Range("A11").Value = 22557
$ws.Range("B11").Value = "DOPE"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 9.784000000000001
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = 9
$ws.Range("G11").Value = 50
$ws.Range("H11").Value = 25
$ws.Range("I11").Value = 20
$ws.Range("J11").Value = 0.02657352221973434
$ws.Range("K11").Value = 2.224227380490118

$ws.Range("A12").Value = 31298
$ws.Range("B12").Value = "DOTAP"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 12.515
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 150
$ws.Range("H12").Value = 45
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 0.9382282290290539
$ws.Range("K12").Value = 9.916769796547158

$ws.Range("A13").Value = 53928
$ws.Range("B13").Value = "DSPC"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 6.12
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 20
$ws.Range("H13").Value = 25
$ws.Range("I13").Value = 40
$ws.Range("J13").Value = 0.02269806244519374
$ws.Range("K13").Value = 2.191526250912545

